$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: was DAIRO CARMONA FRANCESCHI / period 1910 -> now DILSON ALFONSO LOPEZ VALDELAMAR / period 1705
$ws.Range("C18").Value = "1047450132"
$ws.Range("D18").Value = "DILSON ALFONSO LOPEZ VALDELAMAR"
$ws.Range("E18").Value = "1705"
$ws.Range("F18").Value = 9836
$ws.Range("G18").Value = 1475434

# Rows 19-23: same worker (DAIRO CARMONA FRANCESCHI), periods reordered ascending
$ws.Range("E19").Value = "1905"
$ws.Range("E20").Value = "1906"
$ws.Range("E21").Value = "1907"
$ws.Range("E22").Value = "1908"
$ws.Range("E23").Value = "1909"

# Row 24: was DILSON ALFONSO LOPEZ VALDELAMAR / period 1705 -> now DAIRO CARMONA FRANCESCHI / period 1910
$ws.Range("C24").Value = "1047373088"
$ws.Range("D24").Value = "DAIRO CARMONA FRANCESCHI"
$ws.Range("E24").Value = "1910"
$ws.Range("F24").Value = 23187
$ws.Range("G24").Value = 828116
